$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few new prices keep a trailing zero (e.g. "614.30") which Excel's
# automatic number detection would otherwise strip (-> 614.3). Force
# those specific cells to Text so the literal string round-trips.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"

$ws.Range("D2").Value = "69.814.06"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "3.762.50"
$ws.Range("E3").Value = "  +3.81%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "614.30"
$ws.Range("E5").Value = "  +3.51%  "

$ws.Range("D6").Value = "178.23"
$ws.Range("E6").Value = "  -3.47%  "

$ws.Range("D7").Value = "3.761.58"
$ws.Range("E7").Value = "  +3.74%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +4.41%  "

$ws.Range("D11").Value = "6.34"
$ws.Range("E11").Value = "  -2.60%  "

$ws.Range("D12").Value = "0.492"
$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("D13").Value = "40.93"
$ws.Range("E13").Value = "  +5.01%  "

$ws.Range("D14").Value = "0.0000254"
$ws.Range("E14").Value = "  +0.77%  "

$ws.Range("D15").Value = "4.387.13"
$ws.Range("E15").Value = "  +3.81%  "

$ws.Range("D16").Value = "3.762.08"
$ws.Range("E16").Value = "  +3.78%  "

$ws.Range("D17").Value = "69.850.21"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").Value = "514.59"
$ws.Range("E20").Value = "  +1.68%  "

$ws.Range("D21").Value = "16.65"
$ws.Range("E21").Value = "  -2.54%  "

$ws.Range("D22").Value = "9.60"
$ws.Range("E22").Value = "  +3.14%  "

$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").Value = "2.52"
$ws.Range("E24").Value = "  +5.51%  "

$ws.Range("D25").Value = "88.07"
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("D26").Value = "13.32"
$ws.Range("E26").Value = "  -1.00%  "

$ws.Range("D27").Value = "11.13"
$ws.Range("E27").Value = "  +2.76%  "

$ws.Range("D28").Value = "0.0000133"
$ws.Range("E28").Value = "  +21.97%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("D30").Value = "2.49"
$ws.Range("E30").Value = "  -1.52%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.83"
$ws.Range("E31").Value = "  +3.00%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.82"
$ws.Range("E32").Value = "  -3.48%  "

$ws.Range("D33").Value = "31.57"
$ws.Range("E33").Value = "  -3.75%  "

$ws.Range("E34").Value = "  -1.40%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").Value = "6.21"
$ws.Range("E36").Value = "  +1.28%  "

$ws.Range("E37").Value = "  +1.82%  "

$ws.Range("D38").Value = "0.340"
$ws.Range("E38").Value = "  +1.80%  "

$ws.Range("D39").Value = "2.17"
$ws.Range("E39").Value = "  +3.04%  "

$ws.Range("E40").Value = "  +3.15%  "

$ws.Range("D41").Value = "51.24"
$ws.Range("E41").Value = "  +1.22%  "

$ws.Range("D42").Value = "44.65"
$ws.Range("E42").Value = "  -4.20%  "

$ws.Range("D43").Value = "8.79"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "423.08"
$ws.Range("E44").Value = "  +4.50%  "

$ws.Range("D45").Value = "3.062.13"
$ws.Range("E45").Value = "  -1.99%  "

$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -2.26%  "

$ws.Range("D47").Value = "0.0363"
$ws.Range("E47").Value = "  -0.71%  "

$ws.Range("D48").Value = "27.81"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("D49").Value = "2.51"
$ws.Range("E49").Value = "  +2.74%  "

$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").Value = "134.99"
$ws.Range("E51").Value = "  -0.84%  "
